$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Edit 1: first paragraph.
#   "This is a Microsoft word document." ->
#   "This is a Microsoft word document.  (This is a change – Version for
#    main branch)"
# Two trailing spaces are appended to the existing sentence, and the
# parenthetical remark is appended in red, built up across three
# insertions (mirroring the run layout recorded in the target document).
# ----------------------------------------------------------------------

$p1 = $d.Paragraphs(1)
$r = $p1.Range
# Paragraph.Range includes the trailing paragraph mark; trim it off so we
# only touch the visible text.
$textRange = $d.Range($r.Start, $r.End - 1)
$textRange.InsertAfter("  ")
$afterLen = $d.Paragraphs(1).Range.End - 1

$part1 = "(This is a change " + [char]0x2013 + " Ve"
$part2 = "rsion for main branch"
$part3 = ")"

$insertPoint1 = $d.Range($afterLen, $afterLen)
$insertPoint1.InsertAfter($part1)
$p1End = $afterLen + $part1.Length
$range1 = $d.Range($afterLen, $p1End)
$range1.Font.Color = 255

$insertPoint2 = $d.Range($p1End, $p1End)
$insertPoint2.InsertAfter($part2)
$p2End = $p1End + $part2.Length
$range2 = $d.Range($p1End, $p2End)
$range2.Font.Color = 255

$insertPoint3 = $d.Range($p2End, $p2End)
$insertPoint3.InsertAfter($part3)
$p3End = $p2End + $part3.Length
$range3 = $d.Range($p2End, $p3End)
$range3.Font.Color = 255

# ----------------------------------------------------------------------
# Edit 2: drop the trailing "ank God almighty, we are free at last."
# paragraph at the very end of the document (the closing NormalWeb
# paragraph following "Shall be lifted—nevermore!").
# ----------------------------------------------------------------------

$lastIndex = $d.Paragraphs.Count
$d.Paragraphs($lastIndex).Range.Delete()

# ----------------------------------------------------------------------
# Edit 3: the removed paragraph was the last user of several styles, so
# the saved package no longer carries them. Prune the same now-unused
# styles here. Deletions are applied from the highest collection index
# down to the lowest so that earlier indices never shift underneath us
# while we are still working through the list.
# ----------------------------------------------------------------------

$namesToRemove = @(
  "heading 2",
  "heading 4",
  "apple-converted-space",
  "Hyperlink",
  "Heading 2 Char",
  "Heading 4 Char",
  "audio-tool",
  "subscribe",
  "subscribe-more-info",
  "generic-title",
  "podcast-tools__subscribe-links"
)

$styleCount = $d.Styles.Count
$idxToRemove = @()
for ($i = 1; $i -le $styleCount; $i++) {
  $nm = $d.Styles($i).NameLocal
  if ($namesToRemove -contains $nm) {
    $idxToRemove += $i
  }
}

$sortedIdx = $idxToRemove | Sort-Object -Descending
foreach ($i in $sortedIdx) {
  $d.Styles($i).Delete()
}
